# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in AD1:AF1
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Match the header formatting used by the rest of row 1 (e.g. AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record for every data row (2 through 43)
$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 92  # AD - Wins
    $ws.Cells.Item($r, 31).Value2 = 70  # AE - Losses
    $ws.Cells.Item($r, 32).Value2 = 0   # AF - Ties
}
